$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 corresponds to the second data row (file index 1) in the sheet.
$ws.Range("H3").Value  = 0.9111793472080064
$ws.Range("I3").Value  = 0.008115920217427871
$ws.Range("K3").Value  = 123.0192307692308

$ws.Range("Q3").Value  = 6
$ws.Range("R3").Value  = 18
$ws.Range("S3").Value  = 43
$ws.Range("T3").Value  = 84
$ws.Range("U3").Value  = 98

$ws.Range("V3").Value  = 8683
$ws.Range("W3").Value  = 8671
$ws.Range("X3").Value  = 8646
$ws.Range("Y3").Value  = 8605
$ws.Range("Z3").Value  = 8591

$ws.Range("AF3").Value = 0.999309
$ws.Range("AG3").Value = 0.997928
$ws.Range("AH3").Value = 0.995051
$ws.Range("AI3").Value = 0.990333
$ws.Range("AJ3").Value = 0.988721
